$wb = $excel.ActiveWorkbook

# table4_COPR_at_00s80s_inc_raw
$ws = $wb.Worksheets.Item("table4_COPR_at_00s80s_inc_raw")
$ws.Range("B6").Value = 6
$ws.Range("D6").Value = 5
$ws.Range("H6").Value = 5
$ws.Range("B7").Value = 16
$ws.Range("D7").Value = 16
$ws.Range("H7").Value = 16

# table4_COPR_at_00s80s_inc_divto
$ws = $wb.Worksheets.Item("table4_COPR_at_00s80s_inc_divto")
$ws.Range("B2").Value = 0.312
$ws.Range("D2").Value = 0.25
$ws.Range("H2").Value = 0.375
$ws.Range("B3").Value = 0.25
$ws.Range("D3").Value = 0.312
$ws.Range("H4").Value = 0.312
$ws.Range("B5").Value = 0.062
$ws.Range("D5").Value = 0.125
$ws.Range("B6").Value = 0.375
$ws.Range("D6").Value = 0.312
$ws.Range("H6").Value = 0.312
$ws.Range("B7").Value = 16
$ws.Range("D7").Value = 16
$ws.Range("H7").Value = 16

# table4_COPR_at_00s80s_inc_divex
$ws = $wb.Worksheets.Item("table4_COPR_at_00s80s_inc_divex")
$ws.Range("B6").Value = 6
$ws.Range("D6").Value = 5
$ws.Range("H6").Value = 5
$ws.Range("B7").Value = 16
$ws.Range("D7").Value = 16
$ws.Range("H7").Value = 16

# table4_COPR_at_00s80s_dec_raw
$ws = $wb.Worksheets.Item("table4_COPR_at_00s80s_dec_raw")
$ws.Range("B6").Value = 6
$ws.Range("D6").Value = 5
$ws.Range("H6").Value = 5
$ws.Range("B7").Value = 16
$ws.Range("D7").Value = 16
$ws.Range("H7").Value = 16

# table4_COPR_at_00s80s_dec_divto
$ws = $wb.Worksheets.Item("table4_COPR_at_00s80s_dec_divto")
$ws.Range("H3").Value = 0.062
$ws.Range("B5").Value = 0.625
$ws.Range("D5").Value = 0.688
$ws.Range("H5").Value = 0.625
$ws.Range("B6").Value = 0.375
$ws.Range("D6").Value = 0.312
$ws.Range("H6").Value = 0.312
$ws.Range("B7").Value = 16
$ws.Range("D7").Value = 16
$ws.Range("H7").Value = 16

# table4_COPR_at_00s80s_dec_divex
$ws = $wb.Worksheets.Item("table4_COPR_at_00s80s_dec_divex")
$ws.Range("B6").Value = 6
$ws.Range("D6").Value = 5
$ws.Range("H6").Value = 5
$ws.Range("B7").Value = 16
$ws.Range("D7").Value = 16
$ws.Range("H7").Value = 16
